$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/unstyled) used to restore formatting after forcing Text format
$plainStyle = $ws.Range("B2").Style

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.398.77"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.841.68"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "239.30"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "0.6263"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.07440"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "25.00"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").Value = "0.2895"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "0.07718"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.841.14"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "4.978"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "0.6762"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "0.00001030"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "81.84"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "6.233"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "29.433.07"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "232.48"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").Value = "12.32"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "7.327"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "158.01"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "8.488"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "0.1351"
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").Value = "17.34"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "0.07151"
$ws.Range("E28").Value = "  +11.57%  "
$ws.Range("D29").Value = "1.471"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("D30").Value = "1.482"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "4.043"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").Value = "4.033"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").Value = "1.822"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "1.140"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "0.6987"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "2.577"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "0.01842"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "6.918"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("D39").Value = "2.818"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Value = "1.234.47"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "0.9628"
$ws.Range("E41").Value = "  +5.64%  "
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "2.014.63"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "100.98"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "65.53"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("D47").Value = "1.726"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").Value = "6.967"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").Value = "8.932"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("D50").Value = "0.1137"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").Value = "0.3903"
$ws.Range("E51").Value = "  -1.51%  "

$ws.Range("D4").Style = $plainStyle
$ws.Range("D5").Style = $plainStyle
$ws.Range("D6").Style = $plainStyle
$ws.Range("D7").Style = $plainStyle
$ws.Range("D8").Style = $plainStyle
$ws.Range("D9").Style = $plainStyle
$ws.Range("D10").Style = $plainStyle
$ws.Range("D11").Style = $plainStyle
$ws.Range("D13").Style = $plainStyle
$ws.Range("D14").Style = $plainStyle
$ws.Range("D15").Style = $plainStyle
$ws.Range("D16").Style = $plainStyle
$ws.Range("D17").Style = $plainStyle
$ws.Range("D19").Style = $plainStyle
$ws.Range("D20").Style = $plainStyle
$ws.Range("D21").Style = $plainStyle
$ws.Range("D22").Style = $plainStyle
$ws.Range("D23").Style = $plainStyle
$ws.Range("D24").Style = $plainStyle
$ws.Range("D25").Style = $plainStyle
$ws.Range("D26").Style = $plainStyle
$ws.Range("D27").Style = $plainStyle
$ws.Range("D28").Style = $plainStyle
$ws.Range("D29").Style = $plainStyle
$ws.Range("D30").Style = $plainStyle
$ws.Range("D31").Style = $plainStyle
$ws.Range("D32").Style = $plainStyle
$ws.Range("D33").Style = $plainStyle
$ws.Range("D34").Style = $plainStyle
$ws.Range("D35").Style = $plainStyle
$ws.Range("D36").Style = $plainStyle
$ws.Range("D37").Style = $plainStyle
$ws.Range("D38").Style = $plainStyle
$ws.Range("D39").Style = $plainStyle
$ws.Range("D41").Style = $plainStyle
$ws.Range("D42").Style = $plainStyle
$ws.Range("D44").Style = $plainStyle
$ws.Range("D45").Style = $plainStyle
$ws.Range("D46").Style = $plainStyle
$ws.Range("D47").Style = $plainStyle
$ws.Range("D48").Style = $plainStyle
$ws.Range("D49").Style = $plainStyle
$ws.Range("D50").Style = $plainStyle
$ws.Range("D51").Style = $plainStyle
